# Edit script: update Notified_Production_Wind data
# Shifts all timestamps in column A by +5 days (46074.x -> 46079.x)
# and replaces the corresponding production values in column B
# with the latest fetched data (rows 2-97).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aVals = @(46079.01041666666,46079.02083333334,46079.03125,46079.04166666666,46079.05208333334,46079.0625,46079.07291666666,46079.08333333334,46079.09375,46079.10416666666,46079.11458333334,46079.125,46079.13541666666,46079.14583333334,46079.15625,46079.16666666666,46079.17708333334,46079.1875,46079.19791666666,46079.20833333334,46079.21875,46079.22916666666,46079.23958333334,46079.25,46079.26041666666,46079.27083333334,46079.28125,46079.29166666666,46079.30208333334,46079.3125,46079.32291666666,46079.33333333334,46079.34375,46079.35416666666,46079.36458333334,46079.375,46079.38541666666,46079.39583333334,46079.40625,46079.41666666666,46079.42708333334,46079.4375,46079.44791666666,46079.45833333334,46079.46875,46079.47916666666,46079.48958333334,46079.5,46079.51041666666,46079.52083333334,46079.53125,46079.54166666666,46079.55208333334,46079.5625,46079.57291666666,46079.58333333334,46079.59375,46079.60416666666,46079.61458333334,46079.625,46079.63541666666,46079.64583333334,46079.65625,46079.66666666666,46079.67708333334,46079.6875,46079.69791666666,46079.70833333334,46079.71875,46079.72916666666,46079.73958333334,46079.75,46079.76041666666,46079.77083333334,46079.78125,46079.79166666666,46079.80208333334,46079.8125,46079.82291666666,46079.83333333334,46079.84375,46079.85416666666,46079.86458333334,46079.875,46079.88541666666,46079.89583333334,46079.90625,46079.91666666666,46079.92708333334,46079.9375,46079.94791666666,46079.95833333334,46079.96875,46079.97916666666,46079.98958333334,46080)
$bVals = @(920.524,931.064,945.496,949.593,970.5359999999999,997.542,1021.302,1046.905,1077.073,1101.532,1132.267,1156.886,1217.117,1240.529,1264.845,1285.573,1354.253,1377.423,1400.985,1432.609,1502.359,1541.065,1588.608,1623.904,1640.748,1661.487,1686.489,1707.676,1714.817,1709.694,1720.588,1715.12,1668.335,1668.126,1652.964,1649.883,1832.228,1812.022,1787.707,1765.237,1755.211,1722.595,1689.932,1656.253,1625.553,1618.561,1611.877,1605.209,1617.489,1624.588,1631.367,1638.768,1651.801,1659.213,1666.19,1673.011,1709.548,1708.725,1706.367,1703.002,1680.72,1645.357,1609.523,1581.221,1366.736,1310.14,1249.842,1193.69,1120.788,1093.164,1064.168,1038.52,1004.053,989.991,977.901,964.999,976.068,962.5700000000001,951.317,940.407,915.997,903.3,892.0599999999999,879.619,864.264,852.547,841.578,828.5309999999999,805.571,790.02,774.717,759.293,0,0,0,0)

$startRow = 2
$n = $aVals.Length

for ($i = 0; $i -lt $n; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $aVals[$i]
    $ws.Cells.Item($row, 2).Value = $bVals[$i]
}
